$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 11 data rows (rows 2-12), which shifts the remaining
# data (previously rows 13-22) up to rows 2-11.
$ws.Range("A2:C12").EntireRow.Delete()

# The 10 brand-new rows that get appended after the (now shifted) existing
# data, landing in rows 12-21.
$newRows = @(
    @(-0.5496259927749634, 1.477530360221863, -3.358078956604004),
    @(1.193783402442932, 5.348583221435547, -3.039818286895752),
    @(0.0236710291355848, -0.4915938079357147, -0.5149593949317932),
    @(-0.52702397108078, 3.396258115768433, -1.487456917762756),
    @(-0.113315500319004, 3.309820652008057, -0.5609270334243774),
    @(2.999042987823486, 0.8868235945701599, -1.359175205230713),
    @(8.29066276550293, -1.915215253829956, 0.5499314665794373),
    @(2.119396924972534, 0.5285511612892151, 0.0710130855441093),
    @(3.94298243522644, -0.3740022480487823, -1.55419385433197),
    @(3.11648178100586, -1.741576790809631, 3.701537847518921)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
